$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final occupation -> average debt table (occupations kept in alphabetical
# order; three newly classified occupations - mariner, printer and school
# committee - are now included, which pushes several later rows down).
$data = @(
    @("adminastrator", 4081.06),
    @("apothecary", 7.72),
    @("attorney", 11114.6),
    @("blacksmith", 724.7099999999999),
    @("broker", 208.55),
    @("carpenter", 175.34),
    @("cashier", 85.3),
    @("clothier", 53.33),
    @("committee", 863.41),
    @("cooper", 124.76),
    @("doctor", 647.09),
    @("esquire", 13020.09),
    @("executor", 4558.25),
    @("farmer", 36715.82),
    @("gentleman", 16624.32),
    @("goldsmith", 6.82),
    @("guardian", 1245.62),
    @("hatter", 242.99),
    @("inn keeper", 201.53),
    @("joiner", 53.33),
    @("leather dresser", 666.67),
    @("mariner", 3831.05),
    @("merchant", 63765.00000000001),
    @("miller", 18.99),
    @("molster", 211.38),
    @("physician", 1303.46),
    @("post rider", 676.91),
    @("printer", 204.78),
    @("proprietors", 276.9),
    @("sadler", 3524.74),
    @("school committee", 1071.65),
    @("school master", 2213.75),
    @("sherriffs", 512.45),
    @("shipwright", 533.33),
    @("shoemaker", 1345.46),
    @("soldier", 78.48),
    @("spinster", 1731.58),
    @("tanner", 713.0599999999999),
    @("taylor", 209.81),
    @("trader", 9454.129999999999),
    @("treasurer", 15387.15),
    @("widow", 6488.11),
    @("wife", 63.64),
    @("yeoman", 64.54000000000001)
)

$lastRow = 1 + $data.Count

# The sheet originally only goes down to row 42; make sure the new rows
# (43-45) pick up the same formatting (bold/centered/bordered index column)
# as the rest of column A before filling in their values.
$ws.Range("A2").Copy()
$ws.Range("A2:A$lastRow").PasteSpecial(-4122)  # xlPasteFormats

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $row++
}
